$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 24.02.2022 10:30"

# Swap the current/old price values for row 6 (Shell Olomoucká)
$ws.Range("B6").Value = 38.29
$ws.Range("C6").Value = 37.9

# Update delta (kept as text, with sign flipped from negative to positive)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "+0.39"
$ws.Range("D6").Style = "Normal"

# Update the "old datum" timestamp text for row 6
$ws.Range("E6").Value = "2022-02-24 10:30:20"
